$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @("2091", "380:2", "380:5", "380:7", "381:2", "381:4", "381:6", "381:8")

$row2 = $ws.Range("A2:H2")
$row2.NumberFormat = "@"
$row2.HorizontalAlignment = -4108  # xlCenter

for ($i = 0; $i -lt $values.Length; $i++) {
    $cell = $ws.Cells.Item(2, $i + 1)
    $cell.Value = $values[$i]
}
